$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case connector words (de/del/el/la/los/las/y) in state/municipality names
$ws.Range("B8").Value = "Pabellón De Arteaga"
$ws.Range("B9").Value = "Rincón De Romos"
$ws.Range("B10").Value = "San José De Gracia"
$ws.Range("B14").Value = "Playas De Rosarito"
$ws.Range("B31").Value = "Amatenango De La Frontera"
$ws.Range("B32").Value = "Amatenango Del Valle"
$ws.Range("B35").Value = "Bejucal De Ocampo"
$ws.Range("B37").Value = "Benemérito De Las Américas"
$ws.Range("B42").Value = "Chiapa De Corzo"
$ws.Range("B46").Value = "Comitán De Domínguez"
$ws.Range("B73").Value = "Mazapa De Madero"
$ws.Range("B79").Value = "Ocozocoautla De Espinosa"
$ws.Range("B89").Value = "Salto De Agua"
$ws.Range("B90").Value = "San Cristóbal De Las Casas"
$ws.Range("B123").Value = "Hidalgo Del Parral"
$ws.Range("B131").Value = "San Francisco Del Oro"
$ws.Range("B133").Value = "Valle De Zaragoza"
$ws.Range("B151").Value = "San Juan De Sabinas"
$ws.Range("B161").Value = "Villa De Álvarez"
$ws.Range("A163").Value = "Ciudad De México"
$ws.Range("B167").Value = "Cuajimalpa De Morelos"
$ws.Range("B182").Value = "Coneto De Comonfort"
$ws.Range("B194").Value = "Nombre De Dios"
$ws.Range("B198").Value = "Pánuco De Coronado"
$ws.Range("B202").Value = "San Juan De Guadalupe"
$ws.Range("B203").Value = "San Juan Del Río"
$ws.Range("A211").Value = "Estado De México"
$ws.Range("B211").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B214").Value = "Almoloya De Alquisiras"
$ws.Range("B215").Value = "Almoloya De Juárez"
$ws.Range("B219").Value = "Atizapán De Zaragoza"
$ws.Range("B226").Value = "Coacalco De Berriozábal"
$ws.Range("B232").Value = "Ecatepec De Morelos"
$ws.Range("B238").Value = "Ixtapan De La Sal"
$ws.Range("B239").Value = "Ixtapan Del Oro"
$ws.Range("B252").Value = "Naucalpan De Juárez"
$ws.Range("B262").Value = "San Felipe Del Progreso"
$ws.Range("B264").Value = "San Simón De Guerrero"
$ws.Range("B266").Value = "Soyaniquilpan De Juárez"
$ws.Range("B275").Value = "Tenango Del Valle"
$ws.Range("B285").Value = "Tlalnepantla De Baz"
$ws.Range("B290").Value = "Valle De Bravo"
$ws.Range("B291").Value = "Valle De Chalco Solidaridad"
$ws.Range("B292").Value = "Villa De Allende"
$ws.Range("B293").Value = "Villa Del Carbón"
$ws.Range("B305").Value = "Apaseo El Alto"
$ws.Range("B306").Value = "Apaseo El Grande"
$ws.Range("B314").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B318").Value = "Jaral Del Progreso"
$ws.Range("B326").Value = "Purísima Del Rincón"
$ws.Range("B330").Value = "San Diego De La Unión"
$ws.Range("B332").Value = "San Francisco Del Rincón"
$ws.Range("B334").Value = "San Luis De La Paz"
$ws.Range("B336").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B338").Value = "Silao De La Victoria"
$ws.Range("B343").Value = "Valle De Santiago"
$ws.Range("B349").Value = "Acapulco De Juárez"
$ws.Range("B351").Value = "Ajuchitlán Del Progreso"
$ws.Range("B352").Value = "Alcozauca De Guerrero"
$ws.Range("B355").Value = "Atenango Del Río"
$ws.Range("B357").Value = "Atoyac De Álvarez"
$ws.Range("B358").Value = "Ayutla De Los Libres"
$ws.Range("B361").Value = "Buenavista De Cuéllar"
$ws.Range("B362").Value = "Chilapa De Álvarez"
$ws.Range("B363").Value = "Chilpancingo De Los Bravo"
$ws.Range("B364").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B369").Value = "Coyuca De Benítez"
$ws.Range("B370").Value = "Coyuca De Catalán"
$ws.Range("B374").Value = "Cuetzala Del Progreso"
$ws.Range("B375").Value = "Cutzamala De Pinzón"
$ws.Range("B380").Value = "Huitzuco De Los Figueroa"
$ws.Range("B381").Value = "Iguala De La Independencia"
$ws.Range("B383").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B384").Value = "Zihuatanejo De Azueta"
$ws.Range("B386").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B389").Value = "Mártir De Cuilapan"
$ws.Range("B402").Value = "Taxco De Alarcón"
$ws.Range("B404").Value = "Técpan De Galeana"
$ws.Range("B406").Value = "Tepecoacuilco De Trujano"
$ws.Range("B407").Value = "Tixtla De Guerrero"
$ws.Range("B411").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B412").Value = "Tlapa De Comonfort"
$ws.Range("B424").Value = "Agua Blanca De Iturbide"
$ws.Range("B429").Value = "Atotonilco De Tula"
$ws.Range("B430").Value = "Atotonilco El Grande"
$ws.Range("B436").Value = "Cuautepec De Hinojosa"
$ws.Range("B441").Value = "Huasca De Ocampo"
$ws.Range("B445").Value = "Huejutla De Reyes"
$ws.Range("B448").Value = "Jacala De Ledezma"
$ws.Range("B454").Value = "Mineral Del Chico"
$ws.Range("B455").Value = "Mineral Del Monte"
$ws.Range("B456").Value = "Mixquiahuala De Juárez"
$ws.Range("B457").Value = "Molango De Escamilla"
$ws.Range("B459").Value = "Nopala De Villagrán"
$ws.Range("B460").Value = "Omitlán De Juárez"
$ws.Range("B461").Value = "Pachuca De Soto"
$ws.Range("B464").Value = "Progreso De Obregón"
$ws.Range("B469").Value = "Santiago De Anaya"
$ws.Range("B470").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B474").Value = "Tenango De Doria"
$ws.Range("B476").Value = "Tepehuacán De Guerrero"
$ws.Range("B477").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B478").Value = "Tezontepec De Aldama"
$ws.Range("B486").Value = "Tula De Allende"
$ws.Range("B487").Value = "Tulancingo De Bravo"
$ws.Range("B488").Value = "Villa De Tezontepec"
$ws.Range("B491").Value = "Zacualtipán De Ángeles"
$ws.Range("B496").Value = "Ahualulco De Mercado"
$ws.Range("B502").Value = "Atotonilco El Alto"
$ws.Range("B504").Value = "Autlán De Navarro"
$ws.Range("B517").Value = "Encarnación De Díaz"
$ws.Range("B522").Value = "Huejuquilla El Alto"
$ws.Range("B523").Value = "Ixtlahuacán Del Río"
$ws.Range("B527").Value = "Jilotlán De Los Dolores"
$ws.Range("B529").Value = "Lagos De Moreno"
$ws.Range("B536").Value = "Ojuelos De Jalisco"
$ws.Range("B541").Value = "San Diego De Alejandría"
$ws.Range("B542").Value = "San Juan De Los Lagos"
$ws.Range("B544").Value = "San Martín De Bolaños"
$ws.Range("B546").Value = "San Miguel El Alto"
$ws.Range("B547").Value = "San Sebastián Del Oeste"
$ws.Range("B548").Value = "Santa María De Los Ángeles"
$ws.Range("B549").Value = "Santa María Del Oro"
$ws.Range("B551").Value = "Tamazula De Gordiano"
$ws.Range("B557").Value = "Teocuitatlán De Corona"
$ws.Range("B558").Value = "Tepatitlán De Morelos"
$ws.Range("B561").Value = "Tizapán El Alto"
$ws.Range("B562").Value = "Tlajomulco De Zúñiga"
$ws.Range("B569").Value = "Unión De Tula"
$ws.Range("B570").Value = "Valle De Guadalupe"
$ws.Range("B571").Value = "Valle De Juárez"
$ws.Range("B574").Value = "Yahualica De González Gallo"
$ws.Range("B575").Value = "Zacoalco De Torres"
$ws.Range("B577").Value = "Zapotitlán De Vadillo"
$ws.Range("B578").Value = "Zapotlán Del Rey"
$ws.Range("B579").Value = "Zapotlán El Grande"
$ws.Range("B602").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B604").Value = "Cojumatlán De Régules"
$ws.Range("B667").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B694").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B697").Value = "Puente De Ixtla"
$ws.Range("B702").Value = "Tetela Del Volcán"
$ws.Range("B704").Value = "Tlaltizapán De Zapata"
$ws.Range("B715").Value = "Amatlán De Cañas"
$ws.Range("B718").Value = "Ixtlán Del Río"
$ws.Range("B725").Value = "Santa María Del Oro"
$ws.Range("B742").Value = "Mier Y Noriega"
$ws.Range("B743").Value = "Montemorelos"
$ws.Range("B746").Value = "San Nicolás De Los Garza"
$ws.Range("B750").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B754").Value = "Ayoquezco De Aldama"
$ws.Range("B757").Value = "Chalcatongo De Hidalgo"
$ws.Range("B759").Value = "Coicoyán De Las Flores"
$ws.Range("B762").Value = "Constancia Del Rosario"
$ws.Range("B765").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B766").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B767").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B768").Value = "Huautla De Jiménez"
$ws.Range("B769").Value = "Ixtlán De Juárez"
$ws.Range("B770").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B777").Value = "Magdalena Yodocono De Porfirio Díaz"
$ws.Range("B778").Value = "Mariscala De Juárez"
$ws.Range("B781").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B784").Value = "Nejapa De Madero"
$ws.Range("B785").Value = "Oaxaca De Juárez"
$ws.Range("B786").Value = "Ocotlán De Morelos"
$ws.Range("B787").Value = "Pinotepa De Don Luis"
$ws.Range("B789").Value = "Putla Villa De Guerrero"
$ws.Range("B795").Value = "San Agustín De Las Juntas"
$ws.Range("B808").Value = "San Antonino El Alto"
$ws.Range("B810").Value = "San Antonio De La Cal"
$ws.Range("B815").Value = "San Dionisio Del Mar"
$ws.Range("B831").Value = "San José Del Progreso"
$ws.Range("B841").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B877").Value = "San Miguel Del Puerto"
$ws.Range("B879").Value = "San Miguel El Grande"
$ws.Range("B889").Value = "San Pablo Villa De Mitla"
$ws.Range("B892").Value = "San Pedro El Alto"
$ws.Range("B903").Value = "San Pedro Y San Pablo Tequixtepec"
$ws.Range("B935").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B970").Value = "Santo Domingo De Morelos"
$ws.Range("B982").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B984").Value = "Tataltepec De Valdés"
$ws.Range("B985").Value = "Teotitlán De Flores Magón"
$ws.Range("B986").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B987").Value = "Tlacolula De Matamoros"
$ws.Range("B988").Value = "Totontepec Villa De Morelos"
$ws.Range("B990").Value = "Villa De Etla"
$ws.Range("B991").Value = "Villa De Tamazulápam Del Progreso"
$ws.Range("B992").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B993").Value = "Villa De Zaachila"
$ws.Range("B995").Value = "Villa Sola De Vega"
$ws.Range("B996").Value = "Zapotitlán Del Río"
$ws.Range("B998").Value = "Zimatlán De Álvarez"
$ws.Range("B1012").Value = "Ayotoxco De Guerrero"
$ws.Range("B1015").Value = "Chalchicomula De Sesma"
$ws.Range("B1030").Value = "Cuayuca De Andrade"
$ws.Range("B1031").Value = "Cuetzalan Del Progreso"
$ws.Range("B1042").Value = "Huehuetlán El Grande"
$ws.Range("B1048").Value = "Izúcar De Matamoros"
$ws.Range("B1056").Value = "Los Reyes De Juárez"
$ws.Range("B1063").Value = "Palmar De Bravo"
$ws.Range("B1082").Value = "San Salvador El Seco"
$ws.Range("B1083").Value = "San Salvador El Verde"
$ws.Range("B1091").Value = "Tepanco De López"
$ws.Range("B1094").Value = "Tepexi De Rodríguez"
$ws.Range("B1095").Value = "Tetela De Ocampo"
$ws.Range("B1096").Value = "Teteles De Avila Castillo"
$ws.Range("B1101").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B1120").Value = "Xochitlán De Vicente Suárez"
$ws.Range("B1128").Value = "Amealco De Bonfil"
$ws.Range("B1130").Value = "Cadereyta De Montes"
$ws.Range("B1135").Value = "Jalpan De Serra"
$ws.Range("B1136").Value = "Landa De Matamoros"
$ws.Range("B1139").Value = "Pinal De Amoles"
$ws.Range("B1142").Value = "San Juan Del Río"
$ws.Range("B1153").Value = "Armadillo De Los Infante"
$ws.Range("B1154").Value = "Axtla De Terrazas"
$ws.Range("B1160").Value = "Ciudad Del Maíz"
$ws.Range("B1170").Value = "Mexquitic De Carmona"
$ws.Range("B1175").Value = "San Ciro De Acosta"
$ws.Range("B1181").Value = "Santa María Del Río"
$ws.Range("B1183").Value = "Soledad De Graciano Sánchez"
$ws.Range("B1191").Value = "Tanquián De Escobedo"
$ws.Range("B1195").Value = "Villa De Arista"
$ws.Range("B1196").Value = "Villa De Arriaga"
$ws.Range("B1197").Value = "Villa De Guadalupe"
$ws.Range("B1198").Value = "Villa De La Paz"
$ws.Range("B1199").Value = "Villa De Ramos"
$ws.Range("B1200").Value = "Villa De Reyes"
$ws.Range("B1240").Value = "Jalpa De Méndez"
$ws.Range("B1270").Value = "Soto La Marina"
$ws.Range("B1285").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B1286").Value = "Mazatecochco De José María Morelos"
$ws.Range("B1287").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B1289").Value = "Papalotla De Xicohténcatl"
$ws.Range("B1290").Value = "San Pablo Del Monte"
$ws.Range("B1291").Value = "Sanctórum De Lázaro Cárdenas"
$ws.Range("B1296").Value = "Tetla De La Solidaridad"
$ws.Range("B1312").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B1315").Value = "Amatlán De Los Reyes"
$ws.Range("B1323").Value = "Boca Del Río"
$ws.Range("B1328").Value = "Castillo De Teayo"
$ws.Range("B1330").Value = "Cazones De Herrera"
$ws.Range("B1347").Value = "Cosamaloapan De Carpio"
$ws.Range("B1348").Value = "Cosautlán De Carvajal"
$ws.Range("B1364").Value = "Hueyapan De Ocampo"
$ws.Range("B1365").Value = "Ignacio De La Llave"
$ws.Range("B1368").Value = "Ixhuacán De Los Reyes"
$ws.Range("B1369").Value = "Ixhuatlán De Madero"
$ws.Range("B1370").Value = "Ixhuatlán Del Café"
$ws.Range("B1382").Value = "Juchique De Ferrer"
$ws.Range("B1386").Value = "Las Vigas De Ramírez"
$ws.Range("B1387").Value = "Lerdo De Tejada"
$ws.Range("B1390").Value = "Martínez De La Torre"
$ws.Range("B1393").Value = "Medellín De Bravo"
$ws.Range("B1398").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B1405").Value = "Ozuluama De Mascareñas"
$ws.Range("B1408").Value = "Paso Del Macho"
$ws.Range("B1412").Value = "Poza Rica De Hidalgo"
$ws.Range("B1421").Value = "Sayula De Alemán"
$ws.Range("B1425").Value = "Soledad De Doblado"
$ws.Range("B1430").Value = "Tatahuicapan De Juárez"
$ws.Range("B1449").Value = "Tlacotepec De Mejía"
$ws.Range("B1462").Value = "Vega De Alatorre"
$ws.Range("B1471").Value = "Zontecomatlán De López Y Fuentes"
$ws.Range("B1485").Value = "Concepción Del Oro"
$ws.Range("B1493").Value = "Jiménez Del Teul"
$ws.Range("B1501").Value = "Noria De Ángeles"
$ws.Range("B1511").Value = "Teúl De González Ortega"
$ws.Range("B1512").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1513").Value = "Trinidad García De La Cadena"
$ws.Range("B1515").Value = "Villa De Cos"

# Fix floating point percentage serialization noise
$ws.Range("D166").Value = 0.000926964921695858
$ws.Range("D180").Value = 0.000926964921695858
$ws.Range("D215").Value = 0.0009757525491535348
$ws.Range("D220").Value = 0.000926964921695858
$ws.Range("D257").Value = 0.000926964921695858
$ws.Range("D315").Value = 0.000926964921695858
$ws.Range("D340").Value = 0.000926964921695858
$ws.Range("D423").Value = 0.0009757525491535348
$ws.Range("D486").Value = 0.000926964921695858
$ws.Range("D593").Value = 0.000926964921695858
$ws.Range("D681").Value = 0.000926964921695858
$ws.Range("D758").Value = 0.0009757525491535348
$ws.Range("D1087").Value = 0.000926964921695858
$ws.Range("D1149").Value = 0.0009757525491535348
$ws.Range("D1166").Value = 0.0009757525491535348
$ws.Range("D1231").Value = 0.000926964921695858
$ws.Range("D1263").Value = 0.0009757525491535348
$ws.Range("D1272").Value = 0.0009757525491535348
$ws.Range("D1308").Value = 0.000926964921695858
$ws.Range("D1412").Value = 0.000926964921695858
$ws.Range("D1502").Value = 0.000926964921695858
$ws.Range("D1517").Value = 0.000926964921695858
$ws.Range("D1520").Value = 0.0009757525491535348

# Remove trailing footnote rows (1524:1528) and shrink dimension to A1:D1522
$ws.Range("A1524:D1528").EntireRow.Delete()
